$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 232, shifting existing rows 232:265 down to 233:266.
$ws.Range("A232:T232").EntireRow.Insert()

# Populate the newly inserted row 232 with the new price record.
$ws.Cells.Item(232, 1).Value = 10
$ws.Cells.Item(232, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(232, 3).Value = "La Araucanía"
$ws.Cells.Item(232, 4).Value = 45124
$ws.Cells.Item(232, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(232, 5).Value = 9
$ws.Cells.Item(232, 6).Value = "Fruta"
$ws.Cells.Item(232, 7).Value = 100104
$ws.Cells.Item(232, 8).Value = "Frutos de pepita"
$ws.Cells.Item(232, 9).Value = 100104001
$ws.Cells.Item(232, 10).Value = "Granada"
$ws.Cells.Item(232, 11).Value = "Wonderfull"
$ws.Cells.Item(232, 12).Value = "Primera"
$ws.Cells.Item(232, 13).Value = 95
$ws.Cells.Item(232, 14).Value = 13000
$ws.Cells.Item(232, 15).Value = 13000
$ws.Cells.Item(232, 16).Value = 13000
$ws.Cells.Item(232, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(232, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(232, 19).Value = 1300
$ws.Cells.Item(232, 20).Value = 10
